$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false
$ws = $wb.Worksheets.Item("04-15-2022")
$ws.Delete()
